$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Color constants (BGR ints as used by Excel's Interior.Color) ---
$cyan   = 16777164   # FFCCFFFF
$orange = 10079487   # FFFFCC99
$red    = 10066431   # FFFF9999
$green  = 10092390   # FF66FF99
$yellow = 13434879   # FFFFFFCC
$purple = 16751052   # FFCC99FF

# --- 1. Row 1: remove merged "ORIGINAL" header block, replace with plain 0..17 series ---
$ws.Range("H1:K1").UnMerge()
$ws.Range("H1:K1").Style = "Normal"
$headerVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17)
for ($i = 0; $i -lt 18; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerVals[$i]
}

# --- 2. Row 2: updated values ---
$ws.Range("A2").Value = 24
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 3
$ws.Range("O2").Value = 114
$ws.Range("P2").Value = 46

# --- Row 2: recolor in repeating groups of 6 columns (existing borders kept as-is) ---
$ws.Range("B2,H2,N2").Interior.Color = $cyan
$ws.Range("C2,I2,O2").Interior.Color = $orange
$ws.Range("D2,J2,P2").Interior.Color = $red
$ws.Range("E2,K2,Q2").Interior.Color = $green
$ws.Range("F2,L2,R2").Interior.Color = $yellow
$ws.Range("G2,M2").Interior.Color = $purple
$ws.Range("A2").Interior.Color = $purple

# --- 3. New row 3: text cell ---
$ws.Range("A3").Value = "j"

# --- 4. Column widths ---
$ws.Range("A1:R1").ColumnWidth = 4.64

# --- 5. Selection ---
$ws.Range("K3").Select()

# --- 6. Page setup (portrait) ---
$ws.PageSetup.Orientation = 1
